$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header rows: make the repeated header rows bold (reuses existing bold style) ---
$ws.Range("A10:L10").Font.Bold = $true
$ws.Range("A51:M51").Font.Bold = $true

# --- First exchanges table (rows 12-41): drop the stray "uncertainty type = 0" values ---
$ws.Range("K12:K40").ClearContents()

# Row 41's uncertainty flag is actually a boolean TRUE (not a numeric 0)
$ws.Range("K41").Value = $true

# --- Second exchanges table (rows 53-84): same fix, one column over (L instead of K) ---
$ws.Range("L53:L76").ClearContents()
$ws.Range("L81:L83").ClearContents()

# Row 84's uncertainty flag is actually a boolean TRUE (not a numeric 0)
$ws.Range("L84").Value = $true

# --- Restore the view to the top of the sheet with F32 as the active cell ---
$excel.ActiveWindow.ScrollRow = 2
$excel.ActiveWindow.ScrollColumn = 1
[void]$ws.Range("F32").Select()
